# Applies the edits described in the commit "with MAC (after bug fix)".
# This corrects various macProperty.virtue / macProperty.vice columns (AB, AI, BC, BJ)
# for rows 2-17, the derived aggregate row (row 15) formula outputs, and a small
# text correction in the English summary comment (E15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text correction in E15 ---
$ws.Range('E15').Value = 'Sam says: “There is more and more attention for the climate problem in schools. Many young people are concerned about this and that is why it is good that more attention is paid to it. Global warming continues and continues and young people will of course have to live with it the longest. Everyone too easily takes the plane for a week in Bali or the car for a long journey. That really has to change. In addition, it is also very bad to see how animals in the Arctic or in the rainforest lose their habitat because we humans treat the earth badly. The harrowing images of polar bears losing their habitat make a deep impression on me. You can take small measures yourself, such as reducing meat consumption, separating waste and often taking the bike or otherwise the train. If we can convince the new generation, so everyone who is currently in school, that the climate is urgent and needs more attention, that is something very useful. In addition, it would be good if people were given tools to deal with the climate properly. This would create more awareness and better behavior from a large group of people. If these young people then make their parents more aware at home by talking about what they learn about this, that would be great. We have to do this together. So not just a small group, but we all have to show different behavior.”'

# --- Numeric corrections ---
$ws.Range("AB2").Value = 0
$ws.Range("AI2").Value = 0.01164021164021164
$ws.Range("BC2").Value = 0.04119687865012053
$ws.Range("BJ2").Value = 0.0604427577522639
$ws.Range("AB3").Value = 0
$ws.Range("AI3").Value = 0.01042286706349206
$ws.Range("BC3").Value = 0.03538800090408178
$ws.Range("BJ3").Value = 0.06065857317460422
$ws.Range("AB4").Value = 0
$ws.Range("AI4").Value = 0.004025764895330113
$ws.Range("BC4").Value = 0.0291531126509746
$ws.Range("BJ4").Value = 0.07564301370664898
$ws.Range("AB5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("BC5").Value = 0.03974369769837258
$ws.Range("BJ5").Value = 0.04676342710900078
$ws.Range("AB6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("BC6").Value = 0.03662787667868152
$ws.Range("BJ6").Value = 0.04860553693186549
$ws.Range("AB7").Value = 0
$ws.Range("AI7").Value = 0.001484230055658627
$ws.Range("BC7").Value = 0.01993660838121302
$ws.Range("BJ7").Value = 0.07410825350793096
$ws.Range("AB8").Value = 0
$ws.Range("AI8").Value = 0.009564801530368245
$ws.Range("BC8").Value = 0.03797382614171813
$ws.Range("BJ8").Value = 0.07903986214269051
$ws.Range("AB9").Value = 0
$ws.Range("AI9").Value = 0.02930311435970921
$ws.Range("BC9").Value = 0.03460062953910461
$ws.Range("BJ9").Value = 0.07807067059159868
$ws.Range("AB10").Value = 0
$ws.Range("AI10").Value = 0.00582010582010582
$ws.Range("BC10").Value = 0.04360061092080365
$ws.Range("BJ10").Value = 0.05322597395607262
$ws.Range("AB11").Value = 0.003694581280788177
$ws.Range("AI11").Value = 0.003940886699507389
$ws.Range("BC11").Value = 0.0440998709972127
$ws.Range("BJ11").Value = 0.05248791864964339
$ws.Range("AB12").Value = 0.0027124773960217
$ws.Range("AI12").Value = 0.003516174402250352
$ws.Range("BC12").Value = 0.0341636226194029
$ws.Range("BJ12").Value = 0.05910729504110495
$ws.Range("AB13").Value = 0
$ws.Range("AI13").Value = 0.005904796511627907
$ws.Range("BC13").Value = 0.03462978413904382
$ws.Range("BJ13").Value = 0.04404451380011778
$ws.Range("AB14").Value = 0.003456221198156682
$ws.Range("AI14").Value = 0.004480286738351254
$ws.Range("BC14").Value = 0.04103327232924676
$ws.Range("BJ14").Value = 0.0461461324179439
$ws.Range("J15").Value = 0.02029943849002509
$ws.Range("K15").Value = 0.01025076587512379
$ws.Range("L15").Value = 0.006173915710555053
$ws.Range("M15").Value = 0.01001427598489255
$ws.Range("N15").Value = 0.002962962962962963
$ws.Range("O15").Value = 0.0607039732737575
$ws.Range("P15").Value = 0.01262766259525331
$ws.Range("Q15").Value = 0.009251557095840317
$ws.Range("R15").Value = 0.01300802014686134
$ws.Range("S15").Value = 0.006501322751322752
$ws.Range("T15").Value = 3.571428571428572
$ws.Range("U15").Value = 0.0002779651319968244
$ws.Range("V15").Value = 0.03415494617325743
$ws.Range("W15").Value = 0.02001557775853947
$ws.Range("X15").Value = 0.004307692307692308
$ws.Range("Y15").Value = 0.004861019213439956
$ws.Range("Z15").Value = 0.006267806267806268
$ws.Range("AA15").Value = 0.03006969148053615
$ws.Range("AB15").Value = 0
$ws.Range("AC15").Value = 0.04295486169159677
$ws.Range("AD15").Value = 0.003276903276903277
$ws.Range("AE15").Value = 0.01307696421049275
$ws.Range("AF15").Value = 0.008045165843330981
$ws.Range("AG15").Value = 0.04923324490316233
$ws.Range("AH15").Value = 0.01455791914507511
$ws.Range("AI15").Value = 0.005594405594405594
$ws.Range("AJ15").Value = 0.0002553264595400091
$ws.Range("AK15").Value = 0.03530851605774157
$ws.Range("AL15").Value = 0.03839808927112017
$ws.Range("AM15").Value = 0.03726371005668276
$ws.Range("AN15").Value = 0.0357661911482334
$ws.Range("AO15").Value = 0.02524822407770156
$ws.Range("AP15").Value = 0.09709034862728746
$ws.Range("AQ15").Value = 0.06445163665875664
$ws.Range("AR15").Value = 0.05666886477928863
$ws.Range("AS15").Value = 0.06568063883026432
$ws.Range("AT15").Value = 0.06717770305347846
$ws.Range("AU15").Value = 3.571428571428572
$ws.Range("AV15").Value = 0.0004762823239355088
$ws.Range("AW15").Value = 0.07729137002490277
$ws.Range("AX15").Value = 0.09219764906909667
$ws.Range("AY15").Value = 0.05547110082196951
$ws.Range("AZ15").Value = 0.04354799235501969
$ws.Range("BA15").Value = 0.06235645741294775
$ws.Range("BB15").Value = 0.1034763492833253
$ws.Range("BC15").Value = 0.03613942094222568
$ws.Range("BD15").Value = 0.09934866355915074
$ws.Range("BE15").Value = 0.07406413827549599
$ws.Range("BF15").Value = 0.09055461042121636
$ws.Range("BG15").Value = 0.08891573672246431
$ws.Range("BH15").Value = 0.1034072314764283
$ws.Range("BI15").Value = 0.04278344703709514
$ws.Range("BJ15").Value = 0.06206639915133062
$ws.Range("BK15").Value = 0.0005515288936775282
$ws.Range("AB16").Value = 0
$ws.Range("AI16").Value = 0.009078089194429755
$ws.Range("BC16").Value = 0.02989088228105918
$ws.Range("BJ16").Value = 0.08278900721219333
$ws.Range("AB17").Value = 0.003401360544217687
$ws.Range("AI17").Value = 0.02837285676098827
$ws.Range("BC17").Value = 0.04526080992584928
$ws.Range("BJ17").Value = 0.07019632863281694

